$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row above the current row 128. This pushes the
# existing rows 128-239 down to 129-240 (matching the diff's row-shift
# pattern) and automatically grows the sheet dimension to A1:R240.
$ws.Rows.Item(128).Insert()

# Populate the freshly inserted row 128 with the new weekly record.
$ws.Cells.Item(128, 1).Value = 8
$ws.Cells.Item(128, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(128, 3).Value = "Coquimbo"
$ws.Cells.Item(128, 4).Value = 44566
$ws.Cells.Item(128, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(128, 5).Value = 4
$ws.Cells.Item(128, 6).Value = 100112032
$ws.Cells.Item(128, 7).Value = "Zapallo italiano"
$ws.Cells.Item(128, 8).Value = "Sin especificar"
$ws.Cells.Item(128, 9).Value = "Primera"
$ws.Cells.Item(128, 10).Value = 600
$ws.Cells.Item(128, 11).Value = 7000
$ws.Cells.Item(128, 12).Value = 8000
$ws.Cells.Item(128, 13).Value = 7500
$ws.Cells.Item(128, 14).Value = "`$/caja 70 unidades"
$ws.Cells.Item(128, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(128, 16).Value = 107
$ws.Cells.Item(128, 17).Value = 70
$ws.Cells.Item(128, 18).Value = "Hortaliza"
